$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 26321102
$ws.Range("I132").Value = 28576338
$ws.Range("J132").Value = 10002
$ws.Range("K132").Value = 85729014
$ws.Range("L132").Value = 30006
$ws.Range("M132").Value = -85726484
$ws.Range("N132").Value = -35066

$ws.Range("H133").Value = 46666.668
$ws.Range("J133").Value = 46666.668
$ws.Range("L133").Value = 46666.668
$ws.Range("N133").Value = -56786.668

$ws.Range("H134").Value = 48233.375
$ws.Range("J134").Value = 48233.375
$ws.Range("L134").Value = 48233.375
$ws.Range("N134").Value = -58373.375

$ws.Range("H135").Value = 952.4
$ws.Range("I135").Value = 940.5
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 8464.5
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -5929.5
$ws.Range("N135").Value = -14070

$ws.Range("H136").Value = 48184
$ws.Range("J136").Value = 48184
$ws.Range("L136").Value = 48184
$ws.Range("N136").Value = -58384

$ws.Range("H137").Value = 3338.7576
$ws.Range("I137").Value = 1558.1666
$ws.Range("J137").Value = 5475.467
$ws.Range("K137").Value = 4674.4998
$ws.Range("L137").Value = 16426.401
$ws.Range("M137").Value = -2124.4998
$ws.Range("N137").Value = -21526.401

$ws.Range("H138").Value = 5977.172
$ws.Range("I138").Value = 1429.5
$ws.Range("J138").Value = 6270.57
$ws.Range("K138").Value = 4288.5
$ws.Range("L138").Value = 18811.71
$ws.Range("M138").Value = 851.5
$ws.Range("N138").Value = -29091.71

$ws.Range("H139").Value = 40542.195
$ws.Range("J139").Value = 40542.195
$ws.Range("L139").Value = 40542.195
$ws.Range("N139").Value = -50822.195

$ws.Range("H141").Value = 5715.1924
$ws.Range("I141").Value = 5995
$ws.Range("J141").Value = 4176.25
$ws.Range("K141").Value = 17985
$ws.Range("L141").Value = 12528.75
$ws.Range("M141").Value = -12805
$ws.Range("N141").Value = -22888.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 9354.857
$ws.Range("I6").Value = 7500
$ws.Range("J6").Value = 10746
$ws.Range("K6").Value = 7500
$ws.Range("L6").Value = 10746
$ws.Range("M6").Value = -7327
$ws.Range("N6").Value = -11092

$ws.Range("H32").Value = 3339.044
$ws.Range("I32").Value = 3044.923
$ws.Range("K32").Value = 3044.923
$ws.Range("M32").Value = -2757.923

$ws.Range("H61").Value = 1323.8462
$ws.Range("I61").Value = 1008.9
$ws.Range("J61").Value = 2373.6667
$ws.Range("K61").Value = 1008.9
$ws.Range("L61").Value = 2373.6667
$ws.Range("M61").Value = -796.9
$ws.Range("N61").Value = -2797.6667

$ws.Range("H132").Value = 1787.125
$ws.Range("I132").Value = 1044.2979
$ws.Range("J132").Value = 5666.3335
$ws.Range("K132").Value = 3132.8937
$ws.Range("L132").Value = 16999.0005
$ws.Range("M132").Value = -602.8937000000001
$ws.Range("N132").Value = -22059.0005

$ws.Range("H136").Value = 1323.8462
$ws.Range("I136").Value = 1008.9
$ws.Range("J136").Value = 2373.6667
$ws.Range("K136").Value = 3026.7
$ws.Range("L136").Value = 7121.000100000001
$ws.Range("M136").Value = -476.6999999999998
$ws.Range("N136").Value = -12221.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H129").Value = 43489.8
$ws.Range("J129").Value = 43489.8
$ws.Range("L129").Value = 43489.8
$ws.Range("N129").Value = -53489.8

$ws.Range("H134").Value = 1840.4231
$ws.Range("I134").Value = 1107.6833
$ws.Range("J134").Value = 4282.8887
$ws.Range("K134").Value = 3323.0499
$ws.Range("L134").Value = 12848.6661
$ws.Range("M134").Value = -788.0499
$ws.Range("N134").Value = -17918.6661

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5750012.5
$ws.Range("I6").Value = 5750012.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 5750012.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -5749899.5
$ws.Range("N6").ClearContents()

$ws.Range("H31").Value = 7355335
$ws.Range("I31").Value = 1306.8085
$ws.Range("J31").Value = 23814350
$ws.Range("K31").Value = 1306.8085
$ws.Range("L31").Value = 23814350
$ws.Range("M31").Value = -1011.8085
$ws.Range("N31").Value = -23814940

$ws.Range("H34").Value = 7355335
$ws.Range("I34").Value = 1306.8085
$ws.Range("J34").Value = 23814350
$ws.Range("K34").Value = 1306.8085
$ws.Range("L34").Value = 23814350
$ws.Range("M34").Value = -1104.8085
$ws.Range("N34").Value = -23814754

$ws.Range("H58").Value = 1613.1414
$ws.Range("I58").Value = 1532.2113
$ws.Range("J58").Value = 1818.3572
$ws.Range("K58").Value = 1532.2113
$ws.Range("L58").Value = 1818.3572
$ws.Range("M58").Value = -1329.2113
$ws.Range("N58").Value = -2224.3572

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H136").Value = 1613.1414
$ws.Range("I136").Value = 1532.2113
$ws.Range("J136").Value = 1818.3572
$ws.Range("K136").Value = 4596.6339
$ws.Range("L136").Value = 5455.071599999999
$ws.Range("M136").Value = -2046.6339
$ws.Range("N136").Value = -10555.0716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1176.3914
$ws.Range("I5").Value = 355.33334
$ws.Range("J5").Value = 4132.2
$ws.Range("K5").Value = 1066.00002
$ws.Range("L5").Value = 12396.6
$ws.Range("M5").Value = -954.0000199999999
$ws.Range("N5").Value = -12620.6

$ws.Range("H100").Value = 2420
$ws.Range("J100").Value = 2420
$ws.Range("L100").Value = 7260
$ws.Range("N100").Value = -8882

$ws.Range("H114").Value = 3548.7058
$ws.Range("I114").Value = 99.5
$ws.Range("J114").Value = 4008.6
$ws.Range("K114").Value = 298.5
$ws.Range("L114").Value = 12025.8
$ws.Range("M114").Value = 2955.5
$ws.Range("N114").Value = -18533.8

$ws.Range("H131").Value = 845.3134
$ws.Range("I131").Value = 491.1
$ws.Range("J131").Value = 907.4561
$ws.Range("K131").Value = 1473.3
$ws.Range("L131").Value = 2722.3683
$ws.Range("M131").Value = 3566.7
$ws.Range("N131").Value = -12802.3683

$ws.Range("H132").Value = 2377.2856
$ws.Range("J132").Value = 3022.111
$ws.Range("L132").Value = 27198.999
$ws.Range("N132").Value = -32258.999

$ws.Range("H134").Value = 4897.3706
$ws.Range("I134").Value = 4919.9165
$ws.Range("K134").Value = 14759.7495
$ws.Range("M134").Value = -9689.749500000002

$ws.Range("H135").Value = 1176.3914
$ws.Range("I135").Value = 355.33334
$ws.Range("J135").Value = 4132.2
$ws.Range("K135").Value = 3198.00006
$ws.Range("L135").Value = 37189.8
$ws.Range("M135").Value = -663.0000600000003
$ws.Range("N135").Value = -42259.8

$ws.Range("H136").Value = 3424.6667
$ws.Range("I136").Value = 3274.6155
$ws.Range("J136").Value = 4400
$ws.Range("K136").Value = 9823.8465
$ws.Range("L136").Value = 13200
$ws.Range("M136").Value = -4723.8465
$ws.Range("N136").Value = -23400

$ws.Range("H137").Value = 2781.5
$ws.Range("J137").Value = 3650.75
$ws.Range("L137").Value = 10952.25
$ws.Range("N137").Value = -21152.25

$ws.Range("H138").Value = 2145
$ws.Range("I138").Value = 2145
$ws.Range("K138").Value = 6435
$ws.Range("M138").Value = -1295

$ws.Range("H139").Value = 1703.8667
$ws.Range("I139").Value = 1129.8334
$ws.Range("K139").Value = 3389.5002
$ws.Range("M139").Value = 1750.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2883.7827
$ws.Range("I122").Value = 1381.9375
$ws.Range("J122").Value = 6316.5713
$ws.Range("K122").Value = 4145.8125
$ws.Range("L122").Value = 18949.7139
$ws.Range("M122").Value = -1695.8125
$ws.Range("N122").Value = -23849.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4635.4707
$ws.Range("I7").Value = 3218
$ws.Range("J7").Value = 7234.1665
$ws.Range("K7").Value = 3218
$ws.Range("L7").Value = 7234.1665
$ws.Range("M7").Value = -3106
$ws.Range("N7").Value = -7458.1665

$ws.Range("H122").Value = 6194.3125
$ws.Range("I122").Value = 3350.6667
$ws.Range("K122").Value = 10052.0001
$ws.Range("M122").Value = -7602.000100000001

$ws.Range("H126").Value = 4635.4707
$ws.Range("I126").Value = 3218
$ws.Range("J126").Value = 7234.1665
$ws.Range("K126").Value = 9654
$ws.Range("L126").Value = 21702.4995
$ws.Range("M126").Value = -7184
$ws.Range("N126").Value = -26642.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5051783
$ws.Range("I132").Value = 599.3409
$ws.Range("J132").Value = 15154150
$ws.Range("K132").Value = 1798.0227
$ws.Range("L132").Value = 45462450
$ws.Range("M132").Value = 731.9773
$ws.Range("N132").Value = -45467510
